# Insert a new daily price record for Mango (Terminal Hortofrutícola Agro
# Chillán) right above the existing row 84, shifting all subsequent rows
# down by one (old row 84 -> new row 85, ..., old row 143 -> new row 144).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 84; everything at/after 84 shifts down.
$ws.Rows.Item(84).Insert()

# Populate the new row 84 with the new record's data.
$ws.Range("A84").Value = 7
$ws.Range("B84").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C84").Value = 'Ñuble'
$ws.Range("D84").Value = 45090
$ws.Range("E84").Value = 16
$ws.Range("F84").Value = 'Fruta'
$ws.Range("G84").Value = 100108
$ws.Range("H84").Value = 'Tropicales y subtropicales'
$ws.Range("I84").Value = 100108002
$ws.Range("J84").Value = 'Mango'
$ws.Range("K84").Value = 'Sin especificar'
$ws.Range("L84").Value = 'Primera'
$ws.Range("M84").Value = 40
$ws.Range("N84").Value = 10000
$ws.Range("O84").Value = 10000
$ws.Range("P84").Value = 10000
$ws.Range("Q84").Value = '$/bandeja 4 kilos'
$ws.Range("R84").Value = 'Perú'
$ws.Range("S84").Value = 2500
$ws.Range("T84").Value = 4
